# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text in A1 ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.14 = 20386.63 pesos`n✅ 20386.63 pesos = 5.13 = 971.46 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: refresh the scraped rate table ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 194.5
$tasas.Range("O10").Value = 3965.2
$tasas.Range("N12").Value = 3976.96
$tasas.Range("O12").Value = 189.51
